$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row number, new Price (D), new Volume(1h) (E, optional - $null if unchanged)
$updates = @(
    ,@(2, '27.750.11', '  +0.14%  ')
    ,@(3, '1.854.16', '  -0.39%  ')
    ,@(4, '1.017', '  -1.75%  ')
    ,@(5, '320.06', '  -1.06%  ')
    ,@(6, '1.016', '  -1.72%  ')
    ,@(7, '0.4373', '  -0.86%  ')
    ,@(8, '0.3776', '  -0.71%  ')
    ,@(9, '0.07411', '  -0.48%  ')
    ,@(10, '0.8827', '  +0.01%  ')
    ,@(11, '21.54', '  -0.84%  ')
    ,@(12, '1.858.04', '  -0.45%  ')
    ,@(13, '6.774', '  +0.49%  ')
    ,@(14, '5.470', '  -1.47%  ')
    ,@(15, '0.07097', '  -1.54%  ')
    ,@(16, '88.28', '  +5.58%  ')
    ,@(17, '1.021', '  -1.78%  ')
    ,@(18, '0.000008997', '  -0.93%  ')
    ,@(19, '1.016', '  -1.67%  ')
    ,@(20, '15.40', '  -0.79%  ')
    ,@(21, '27.749.40', '  +0.05%  ')
    ,@(22, '5.262', '  -0.80%  ')
    ,@(23, '11.13', '  -2.58%  ')
    ,@(24, '2.085.99', '  -0.01%  ')
    ,@(25, '2.028', '  +5.05%  ')
    ,@(26, '156.69', '  -1.19%  ')
    ,@(27, '18.67', '  -0.75%  ')
    ,@(28, '5.443', '  +2.48%  ')
    ,@(29, '1.984', '  -0.47%  ')
    ,@(30, '120.59', '  +2.65%  ')
    ,@(31, '0.09038', '  -0.54%  ')
    ,@(32, '1.227', '  +1.43%  ')
    ,@(33, '0.7680', '  +0.34%  ')
    ,@(34, '2.996', '  +3.48%  ')
    ,@(35, '4.549', '  -0.40%  ')
    ,@(36, '1.017', '  -1.63%  ')
    ,@(37, '1.135', '  -1.80%  ')
    ,@(38, '0.01977', '  -0.30%  ')
    ,@(39, '0.05300', $null)
    ,@(40, '2.866', '  +1.30%  ')
    ,@(41, '0.5180', '  -0.23%  ')
    ,@(42, '6.946', '  +1.41%  ')
    ,@(43, '0.1675', '  -0.72%  ')
    ,@(44, '8.692', '  +0.34%  ')
    ,@(45, '109.96', '  +0.41%  ')
    ,@(46, '10.70', '  +1.16%  ')
    ,@(47, '1.709', '  -0.71%  ')
    ,@(48, '0.4711', '  +0.67%  ')
    ,@(49, '1.017', '  -1.66%  ')
    ,@(50, '0.06466', '  +0.79%  ')
    ,@(51, '1.849', '  -0.42%  ')
)

# Protect the Price column values (e.g. "1.017", "320.06") from being auto-
# converted into numbers when assigned via .Value, by temporarily formatting
# the column as Text; ClearFormats() afterwards restores the original (no
# explicit style) look of the cells while keeping the text content intact.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($u in $updates) {
    $r = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]

    $ws.Cells.Item($r, 4).Value = $dVal

    if ($null -ne $eVal) {
        $ws.Cells.Item($r, 5).Value = $eVal
    }
}

$priceRange.ClearFormats()
